# Refresh the NATMI ligand/receptor edge-weight statistics for Jam2-Jam3
# with values recomputed from the updated TPM expression matrix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues($rowNum, $values) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$rowNum").Value = [double]$values[$col]
    }
}

Set-RowValues 2 @{ "E"="3"; "F"="1"; "G"="115.5575153333333"; "H"="346.672546"; "I"="0.9048104954928987"; "J"="0.9048104954928987"; "M"="17.78713366666667"; "N"="53.361401"; "O"="0.2123673935064285"; "P"="0.2123673935064285"; "Q"="2055.436971421883"; "R"="18498.93274279695"; "S"="0.192152246545087"; "T"="0.192152246545087" }
Set-RowValues 3 @{ "E"="3"; "F"="1"; "G"="115.5575153333333"; "H"="346.672546"; "I"="0.9048104954928987"; "J"="0.9048104954928987"; "O"="0.1369154545457259"; "P"="0.1369154545457259"; "Q"="1325.161469403251"; "R"="11926.45322462926"; "S"="0.1238825402681537"; "T"="0.1238825402681537" }
Set-RowValues 4 @{ "E"="3"; "F"="1"; "G"="115.5575153333333"; "H"="346.672546"; "I"="0.9048104954928987"; "J"="0.9048104954928987"; "M"="53.74594866666666"; "N"="161.237846"; "O"="0.6416934422244821"; "P"="0.6416934422244821"; "Q"="6210.74828715288"; "R"="55896.73458437592"; "S"="0.5806109614136774"; "T"="0.5806109614136774" }
Set-RowValues 5 @{ "E"="3"; "F"="1"; "G"="115.5575153333333"; "H"="346.672546"; "I"="0.9048104954928987"; "J"="0.9048104954928987"; "M"="0.7557936666666666"; "N"="2.267381"; "O"="0.009023709723363511"; "P"="0.009023709723363511"; "Q"="87.33763822466955"; "R"="786.038744022026"; "S"="0.008164747265980626"; "T"="0.008164747265980626" }
Set-RowValues 6 @{ "G"="5.519651666666666"; "I"="0.04321864090845719"; "J"="0.04321864090845719"; "M"="17.78713366666667"; "N"="53.361401"; "O"="0.2123673935064285"; "P"="0.2123673935064285"; "Q"="98.17878198843943"; "R"="883.6090378959549"; "S"="0.009178230120619359"; "T"="0.009178230120619359" }
Set-RowValues 7 @{ "G"="5.519651666666666"; "I"="0.04321864090845719"; "J"="0.04321864090845719"; "O"="0.1369154545457259"; "P"="0.1369154545457259"; "Q"="63.29687595043166"; "R"="569.671883553885"; "S"="0.005917299864829921"; "T"="0.005917299864829921" }
Set-RowValues 8 @{ "G"="5.519651666666666"; "I"="0.04321864090845719"; "J"="0.04321864090845719"; "M"="53.74594866666666"; "N"="161.237846"; "O"="0.6416934422244821"; "P"="0.6416934422244821"; "Q"="296.6589151345477"; "R"="2669.930236210929"; "S"="0.02773311845281172"; "T"="0.02773311845281172" }
Set-RowValues 9 @{ "G"="5.519651666666666"; "I"="0.04321864090845719"; "J"="0.04321864090845719"; "K"="2"; "L"="0.6666666666666666"; "M"="0.7557936666666666"; "N"="2.267381"; "O"="0.009023709723363511"; "P"="0.009023709723363511"; "Q"="4.171717771872776"; "R"="37.54545994685499"; "S"="0.0003899924701962012"; "T"="0.0003899924701962012" }
Set-RowValues 10 @{ "G"="6.580297333333334"; "H"="19.740892"; "I"="0.05152345196666309"; "J"="0.05152345196666309"; "M"="17.78713366666667"; "N"="53.361401"; "O"="0.2123673935064285"; "P"="0.2123673935064285"; "Q"="117.0446282344102"; "R"="1053.401654109692"; "S"="0.01094190119861391"; "T"="0.01094190119861391" }
Set-RowValues 11 @{ "G"="6.580297333333334"; "H"="19.740892"; "I"="0.05152345196666309"; "J"="0.05152345196666309"; "O"="0.1369154545457259"; "P"="0.1369154545457259"; "Q"="75.45988210456935"; "R"="679.1389389411241"; "S"="0.007054356845780552"; "T"="0.007054356845780552" }
Set-RowValues 12 @{ "G"="6.580297333333334"; "H"="19.740892"; "I"="0.05152345196666309"; "J"="0.05152345196666309"; "M"="53.74594866666666"; "N"="161.237846"; "O"="0.6416934422244821"; "P"="0.6416934422244821"; "Q"="353.6643226887369"; "R"="3182.978904198632"; "S"="0.0330622612477758"; "T"="0.0330622612477758" }
Set-RowValues 13 @{ "G"="6.580297333333334"; "H"="19.740892"; "I"="0.05152345196666309"; "J"="0.05152345196666309"; "K"="2"; "L"="0.6666666666666666"; "M"="0.7557936666666666"; "N"="2.267381"; "O"="0.009023709723363511"; "P"="0.009023709723363511"; "Q"="4.973347049316889"; "R"="44.760123443852"; "S"="0.0004649326744928305"; "T"="0.0004649326744928305" }
Set-RowValues 14 @{ "E"="1"; "F"="0.3333333333333333"; "G"="0.057141"; "H"="0.171423"; "I"="0.0004474116319810314"; "J"="0.0004474116319810314"; "M"="17.78713366666667"; "N"="53.361401"; "O"="0.2123673935064285"; "P"="0.2123673935064285"; "Q"="1.016374604847"; "R"="9.147371443622999"; "S"="9.501564210826907E-05"; "T"="9.501564210826907E-05" }
Set-RowValues 15 @{ "E"="1"; "F"="0.3333333333333333"; "G"="0.057141"; "H"="0.171423"; "I"="0.0004474116319810314"; "J"="0.0004474116319810314"; "O"="0.1369154545457259"; "P"="0.1369154545457259"; "Q"="0.655267217409"; "R"="5.897404956681"; "S"="6.125756696172795E-05"; "T"="6.125756696172795E-05" }
Set-RowValues 16 @{ "E"="1"; "F"="0.3333333333333333"; "G"="0.057141"; "H"="0.171423"; "I"="0.0004474116319810314"; "J"="0.0004474116319810314"; "M"="53.74594866666666"; "N"="161.237846"; "O"="0.6416934422244821"; "P"="0.6416934422244821"; "Q"="3.071097252762"; "R"="27.639875274858"; "S"="0.0002871011102171812"; "T"="0.0002871011102171812" }
Set-RowValues 17 @{ "E"="1"; "F"="0.3333333333333333"; "G"="0.057141"; "H"="0.171423"; "I"="0.0004474116319810314"; "J"="0.0004474116319810314"; "K"="2"; "L"="0.6666666666666666"; "M"="0.7557936666666666"; "N"="2.267381"; "O"="0.009023709723363511"; "P"="0.009023709723363511"; "Q"="0.04318680590699999"; "R"="0.388681253163"; "S"="0.0004037312693853169"; "T"="0.0004037312693853169" }
